# Applies the cryptos.xlsx price/volume refresh described in the commit diff
# (coin rankings list updated by the scheduled GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns D/E hold numeric-looking text (e.g. '260.30', '0.619', '-1.90%')
# that must stay literal strings rather than being coerced to floating point
# numbers (which would lose trailing zeros / exact decimal text and reformat
# thousand-grouped values). Force Text format on any cell whose new value
# would otherwise parse as a plain number before writing it.

$ws.Range('D2').Value = '37.447.22'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '2.015.85'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.30'
$ws.Range('E5').Value = '  +4.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.619'
$ws.Range('E6').Value = '  -1.90%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.13'
$ws.Range('E8').Value = '  -5.65%  '
$ws.Range('E9').Value = '  -3.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0777'
$ws.Range('E10').Value = '  -4.17%  '
$ws.Range('E11').Value = '  -3.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.35'
$ws.Range('E12').Value = '  -6.28%  '
$ws.Range('D13').Value = '2.310.33'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.57'
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.798'
$ws.Range('E15').Value = '  -7.72%  '
$ws.Range('E16').Value = '  -5.26%  '
$ws.Range('D17').Value = '2.026.80'
$ws.Range('D18').Value = '37.270.80'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.99'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').Value = '0.0₃0838'
$ws.Range('E20').Value = '  -3.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '232.71'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').Value = '  -2.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.60'
$ws.Range('E23').Value = '  +3.62%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.14'
$ws.Range('E26').Value = '  +0.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.98'
$ws.Range('E27').Value = '  -5.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.61'
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('E29').Value = '  -5.82%  '
$ws.Range('E30').Value = '  -4.95%  '
$ws.Range('E31').Value = '  -2.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.61'
$ws.Range('E32').Value = '  -4.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0644'
$ws.Range('E33').Value = '  -3.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.54'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('E35').Value = '  -5.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.81'
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  -4.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.50'
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('E40').Value = '  +3.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.20'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0212'
$ws.Range('E42').Value = '  -1.36%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0931'
$ws.Range('E43').Value = '  -5.27%  '
$ws.Range('D44').Value = '1.432.52'
$ws.Range('E44').Value = '  +3.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.81'
$ws.Range('E45').Value = '  -8.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.80'
$ws.Range('E46').Value = '  -3.20%  '
$ws.Range('E47').Value = '  -3.14%  '
$ws.Range('E48').Value = '  +2.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.03'
$ws.Range('E49').Value = '  -6.79%  '
$ws.Range('D50').Value = '2.201.33'
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.96'
$ws.Range('E51').Value = '  -10.17%  '
